$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "ATATTAATAT"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "'5"
$ws.Range("H3").ClearFormats()

$ws.Range("B4").Value = "ATATATAAAT"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = "'7"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = ""

$ws.Range("B5").Value = "ATATATTAAT"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = "'7"
$ws.Range("H5").ClearFormats()

$ws.Range("B6").Value = "TAATATATAT"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "'1"
$ws.Range("H6").ClearFormats()

$ws.Range("B7").Value = "ATAAATATAT"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = "'3"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").Value = ""

$ws.Range("B8").Value = "AAATATATAT"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("G8").Value = "'1"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").Value = ""

$ws.Range("B9").Value = "ATATAAATAT"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = "'5"
$ws.Range("G9").ClearFormats()
$ws.Range("H9").Value = ""

$ws.Range("B11").Value = "AAATATAAAT"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0
$ws.Range("G11").Value = "1, 7"
$ws.Range("H11").Value = ""

$ws.Range("B12").Value = "AAAAATATAT"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = "1, 3"
$ws.Range("H12").Value = ""

$ws.Range("B13").Value = "ATAAATAAAT"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = "3, 7"
$ws.Range("H13").Value = ""

$ws.Range("B15").Value = "ATTAATAAAT"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("G15").Value = "'7"
$ws.Range("G15").ClearFormats()
$ws.Range("H15").Value = "'3"
$ws.Range("H15").ClearFormats()

$ws.Range("B16").Value = "AATAATATAT"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = "'1"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").Value = "'3"
$ws.Range("H16").ClearFormats()

$ws.Range("B17").Value = "ATAAAAATAT"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = "3, 5"
$ws.Range("H17").Value = ""

$ws.Range("B18").Value = "AAATAAATAT"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = "1, 5"
$ws.Range("H18").Value = ""

$ws.Range("B19").Value = "AAATATTAAT"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("G19").Value = "'1"
$ws.Range("G19").ClearFormats()
$ws.Range("H19").Value = "'7"
$ws.Range("H19").ClearFormats()

$ws.Range("B20").Value = "ATATAAAAAT"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = "5, 7"
$ws.Range("H20").Value = ""

$ws.Range("B21").Value = "TAATATAAAT"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("G21").Value = "'7"
$ws.Range("G21").ClearFormats()
$ws.Range("H21").Value = "'1"
$ws.Range("H21").ClearFormats()

$ws.Range("B23").Value = "ATAATAATAT"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("G23").Value = "'3"
$ws.Range("G23").ClearFormats()
$ws.Range("H23").Value = "'5"
$ws.Range("H23").ClearFormats()

$ws.Range("B24").Value = "TAAAATATAT"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("G24").Value = "'3"
$ws.Range("G24").ClearFormats()
$ws.Range("H24").Value = "'1"
$ws.Range("H24").ClearFormats()

$ws.Range("B25").Value = "AAATAATAAT"
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 1
$ws.Range("G25").Value = "1, 5"
$ws.Range("H25").Value = "'7"
$ws.Range("H25").ClearFormats()

$ws.Range("B26").Value = "ATAAAAAAAT"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = "3, 5, 7"
$ws.Range("H26").Value = ""

$ws.Range("B27").Value = "AAAAATAAAT"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = "1, 3, 7"
$ws.Range("H27").Value = ""

$ws.Range("B28").Value = "AATAATAAAT"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 1
$ws.Range("G28").Value = "1, 7"
$ws.Range("H28").Value = "'3"
$ws.Range("H28").ClearFormats()

$ws.Range("B30").Value = "AAAAAAATAT"
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 0
$ws.Range("G30").Value = "1, 3, 5"
$ws.Range("H30").Value = ""

$ws.Range("B31").Value = "AAAATAATAT"
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 1
$ws.Range("G31").Value = "1, 3"
$ws.Range("H31").Value = "'5"
$ws.Range("H31").ClearFormats()

$ws.Range("B32").Value = "TAAAATAAAT"
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 1
$ws.Range("G32").Value = "3, 7"
$ws.Range("H32").Value = "'1"
$ws.Range("H32").ClearFormats()
